# Update the "cryptos" price/volume table (GitHub Actions daily refresh).
# Price (col D) and Volume(1h) (col E) are refreshed for every coin row;
# rows 41/42 additionally swap their Coin/Link/Price/Volume (TheSandbox <-> Frax).
# D values that look like plain numbers ("330.94", "1.002", ...) are forced to
# Text format first so Excel doesn't silently coerce them into real numbers
# (which would also drop meaningful trailing zeros, e.g. "5.950" -> 5.95).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.640.64'
$ws.Range("E2").Value = '  -1.60%  '

# Row 3
$ws.Range("D3").Value = '1.878.65'
$ws.Range("E3").Value = '  -1.29%  '

# Row 4
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '330.94'
$ws.Range("E5").Value = '  +1.13%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  -0.02%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4711'
$ws.Range("E7").Value = '  +1.88%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3965'
$ws.Range("E8").Value = '  -0.07%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.96'
$ws.Range("E9").Value = '  -7.86%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08057'
$ws.Range("E10").Value = '  -3.02%  '

# Row 11
$ws.Range("E11").Value = '  -1.65%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.78'
$ws.Range("E12").Value = '  -0.32%  '

# Row 13
$ws.Range("D13").Value = '1.884.96'
$ws.Range("E13").Value = '  -0.74%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.950'
$ws.Range("E14").Value = '  -1.12%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.184'
$ws.Range("E15").Value = '  -2.38%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.004'
$ws.Range("E16").Value = '  +0.00%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '86.84'
$ws.Range("E17").Value = '  -2.82%  '

# Row 18
$ws.Range("E18").Value = '  -2.36%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06579'
$ws.Range("E19").Value = '  -0.07%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.25'
$ws.Range("E20").Value = '  -3.03%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  -0.18%  '

# Row 22
$ws.Range("D22").Value = '27.663.31'
$ws.Range("E22").Value = '  -1.44%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.501'
$ws.Range("E23").Value = '  -3.39%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.98'
$ws.Range("E24").Value = '  -1.60%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.308'
$ws.Range("E25").Value = '  -0.35%  '

# Row 26
$ws.Range("D26").Value = '2.104.72'
$ws.Range("E26").Value = '  -1.02%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '154.56'
$ws.Range("E27").Value = '  +0.21%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.24'
$ws.Range("E28").Value = '  +1.16%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.095'
$ws.Range("E29").Value = '  -1.33%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.587'
$ws.Range("E30").Value = '  -2.19%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '122.39'
$ws.Range("E31").Value = '  -1.87%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09493'
$ws.Range("E32").Value = '  -1.05%  '

# Row 33
$ws.Range("E33").Value = '  -1.27%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.470'
$ws.Range("E34").Value = '  +0.13%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.617'
$ws.Range("E35").Value = '  -0.21%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.299'
$ws.Range("E36").Value = '  -3.71%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06108'
$ws.Range("E37").Value = '  -0.55%  '

# Row 38
$ws.Range("E38").Value = '  -1.59%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.220'
$ws.Range("E39").Value = '  -3.30%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.212'
$ws.Range("E40").Value = '  -5.54%  '

# Row 41
$ws.Range("B41").Value = 'Frax'
$ws.Range("C41").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.001'
$ws.Range("E41").Value = '  -0.05%  '

# Row 42
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5992'
$ws.Range("E42").Value = '  -2.11%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1896'
$ws.Range("E43").Value = '  -0.29%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.32'
$ws.Range("E44").Value = '  -4.88%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5697'

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.249'
$ws.Range("E46").Value = '  -3.85%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.15'
$ws.Range("E47").Value = '  -4.64%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.413'
$ws.Range("E48").Value = '  -0.54%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.940'
$ws.Range("E49").Value = '  -3.27%  '

# Row 50
$ws.Range("E50").Value = '  -1.09%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '110.02'
$ws.Range("E51").Value = '  -0.65%  '
